$wb = $excel.ActiveWorkbook

$wsGen  = $wb.Worksheets.Item("GeneralVariables")

# --- GeneralVariables: insert the two new "testingCompanySOI7xx" rows
# right after the existing "testingCompanySOI718" row (originally row 11,
# i.e. before the old row 12 "idTestingCompanySOI66").
$wsGen.Rows.Item(12).Insert()
$wsGen.Rows.Item(13).Insert()
$wsGen.Range("A12").Value = "testingCompanySOI720"
$wsGen.Range("A13").Value = "testingCompanySOI770"
$wsGen.Range("B12").Value = "AutoTestingCompany_SOI720"
$wsGen.Range("B13").Value = "AutoTestingCompany_SOI770"

# --- GeneralVariables: insert the two new "idTestingCompanySOI7xx" rows
# right after the existing "idTestingCompanySOI718" row (now row 23 after
# the two inserts above), before the "optyStage" row.
$wsGen.Rows.Item(24).Insert()
$wsGen.Rows.Item(25).Insert()
$wsGen.Range("B25").Value = "0013E00001AAevfQAD"
$wsGen.Range("A24").Value = "idTestingCompanySOI720"
$wsGen.Range("A25").Value = "idTestingCompanySOI770"
$wsGen.Range("B24").Value = "0013E00001AAevVQAT"

# --- Tab / selection bookkeeping: GeneralVariables becomes the active
# sheet/tab (with a new selection). Activating it implicitly clears the
# previously-selected tab (Environment_DirectSales).
$wsGen.Activate()
[void]$wsGen.Range("I14").Select()
